# Updated cryptos list on Mon Jul 31 17:28:15 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for rows 2-51.
# Some Price values are plain decimals (e.g. "1.001") that Excel would
# otherwise auto-convert to a Number; force those cells to Text first via
# NumberFormat "@" so the literal string is preserved exactly, then restore
# the cell's style to "Normal" so no stray number-format/style is left
# behind (matches the source cells, which carry no explicit style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.192.81'
$ws.Range("E2").Value = '  -0.66%  '
$ws.Range("D3").Value = '1.858.23'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6988'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.57%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07813'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.81%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3115'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.00'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.91%  '
$ws.Range("E11").Value = '  -4.07%  '
$ws.Range("D12").Value = '1.857.09'
$ws.Range("E12").Value = '  -1.39%  '
$ws.Range("E13").Value = '  -1.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6923'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.595'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008519'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.73%  '
$ws.Range("D18").Value = '29.221.80'
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '248.13'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.11%  '
$ws.Range("D20").Value = '2.110.41'
$ws.Range("E20").Value = '  -1.26%  '
$ws.Range("E21").Value = '  -3.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9995'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.570'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.38%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1536'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.74'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.922'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.57'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.572'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.275'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.243'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.31%  '
$ws.Range("E32").Value = '  -1.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05244'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7592'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.871'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.88%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.175'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.702'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01857'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.82%  '
$ws.Range("D39").Value = '1.238.63'
$ws.Range("E39").Value = '  -2.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.738'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9020'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.98'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.874'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -8.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9999'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '68.32'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.00%  '
$ws.Range("D46").Value = '2.010.72'
$ws.Range("E46").Value = '  -1.10%  '
$ws.Range("E47").Value = '  -4.31%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5179'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.508'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.766'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4254'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.13%  '
